# Update Name of Algo
# Applies the data value corrections captured by the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.625
$ws.Range("C3").Value = -12.774
$ws.Range("E19").Value = 16.419
$ws.Range("A21").Value = -20.423
$ws.Range("A23").Value = -20.585
$ws.Range("C24").Value = -12.662
$ws.Range("E24").Value = 17.137
$ws.Range("A25").Value = -21.796
$ws.Range("B27").Value = 5.513
$ws.Range("E30").Value = 16.342
$ws.Range("B31").Value = 5.837999999999999
$ws.Range("E31").Value = 16.255
$ws.Range("E33").Value = 17.229
$ws.Range("B39").Value = 7.997
$ws.Range("B48").Value = 5.274
$ws.Range("B51").Value = 5.970999999999999
$ws.Range("B52").Value = 5.705
$ws.Range("A53").Value = -21.817
$ws.Range("B55").Value = 4.631
$ws.Range("E55").Value = 16.491
$ws.Range("B56").Value = 4.906999999999999
$ws.Range("A57").Value = -21.352
$ws.Range("B57").Value = 5.948
$ws.Range("C57").Value = -13.287
$ws.Range("A59").Value = -22.208
$ws.Range("C61").Value = -13.508
$ws.Range("E65").Value = 17.336
$ws.Range("A69").Value = -21.649
$ws.Range("C70").Value = -11.71
$ws.Range("E70").Value = 17.447
$ws.Range("B73").Value = 7.343000000000001
$ws.Range("E75").Value = 16.682
$ws.Range("A79").Value = -21.192
$ws.Range("A83").Value = -21.943
$ws.Range("E83").Value = 16.605
$ws.Range("C86").Value = -13.597
$ws.Range("B89").Value = 5.986999999999999
$ws.Range("B90").Value = 5.833
$ws.Range("A93").Value = -21.476
$ws.Range("E96").Value = 16.225
$ws.Range("E97").Value = 16.881
$ws.Range("C98").Value = -12.285
$ws.Range("C100").Value = -13.174
$ws.Range("C102").Value = -13.564

$wb.Save()
